$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1440.2
$ws.Range("I40").Value = 1500.5
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1500.5
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -1325.5
$ws.Range("N40").Value = -1750

$ws.Range("H64").Value = 2780.7693
$ws.Range("I64").Value = 2728.5715
$ws.Range("J64").Value = 2841.6667
$ws.Range("K64").Value = 2728.5715
$ws.Range("L64").Value = 2841.6667
$ws.Range("M64").Value = -2480.5715
$ws.Range("N64").Value = -3337.6667

$ws.Range("H67").Value = 2780.7693
$ws.Range("I67").Value = 2728.5715
$ws.Range("J67").Value = 2841.6667
$ws.Range("K67").Value = 2728.5715
$ws.Range("L67").Value = 2841.6667
$ws.Range("M67").Value = -1870.5715
$ws.Range("N67").Value = -4557.6667

$ws.Range("H74").Value = 8339146.5
$ws.Range("I74").Value = 11115556
$ws.Range("J74").Value = 9919.666999999999
$ws.Range("K74").Value = 11115556
$ws.Range("L74").Value = 9919.666999999999
$ws.Range("M74").Value = -11114620
$ws.Range("N74").Value = -11791.667

$ws.Range("H76").Value = 3180.889
$ws.Range("I76").Value = 3168
$ws.Range("J76").Value = 3400
$ws.Range("K76").Value = 3168
$ws.Range("L76").Value = 3400
$ws.Range("M76").Value = -2853
$ws.Range("N76").Value = -4030

$ws.Range("H77").Value = 8339146.5
$ws.Range("I77").Value = 11115556
$ws.Range("J77").Value = 9919.666999999999
$ws.Range("K77").Value = 55577780
$ws.Range("L77").Value = 49598.335
$ws.Range("M77").Value = -55573100
$ws.Range("N77").Value = -58958.335

$ws.Range("H79").Value = 3180.889
$ws.Range("I79").Value = 3168
$ws.Range("J79").Value = 3400
$ws.Range("K79").Value = 3168
$ws.Range("L79").Value = 3400
$ws.Range("M79").Value = -2076
$ws.Range("N79").Value = -5584

$ws.Range("H113").Value = 9308.154
$ws.Range("I113").Value = 4250
$ws.Range("J113").Value = 13643.714
$ws.Range("K113").Value = 4250
$ws.Range("L113").Value = 13643.714
$ws.Range("M113").Value = -996
$ws.Range("N113").Value = -20151.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null

$ws.Range("H63").Value = 6025231
$ws.Range("I63").Value = 10656786
$ws.Range("J63").Value = 4210
$ws.Range("K63").Value = 10656786
$ws.Range("L63").Value = 4210
$ws.Range("M63").Value = -10656100
$ws.Range("N63").Value = -5582

$ws.Range("H66").Value = 6025231
$ws.Range("I66").Value = 10656786
$ws.Range("J66").Value = 4210
$ws.Range("K66").Value = 53283930
$ws.Range("L66").Value = 21050
$ws.Range("M66").Value = -53280498
$ws.Range("N66").Value = -27914

$ws.Range("H130").Value = 43266
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43266
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43266
$ws.Range("M130").Value = $null
$ws.Range("N130").Value = -53306

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 35000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35480

$ws.Range("H27").Value = 35000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 35000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 35000
$ws.Range("N27").Value = -35384

$ws.Range("H123").Value = 38780
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 38780
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 38780
$ws.Range("N123").Value = -48580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8933060
$ws.Range("I131").Value = 125052616
$ws.Range("J131").Value = 786.46155
$ws.Range("K131").Value = 375157848
$ws.Range("L131").Value = 2359.38465
$ws.Range("M131").Value = -375152808
$ws.Range("N131").Value = -12439.38465

$ws.Range("H137").Value = 3930.7058
$ws.Range("I137").Value = 3120.818
$ws.Range("J137").Value = 5415.5
$ws.Range("K137").Value = 9362.454000000002
$ws.Range("L137").Value = 16246.5
$ws.Range("M137").Value = -4262.454000000002
$ws.Range("N137").Value = -26446.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6224.5
$ws.Range("I70").Value = 5624.353
$ws.Range("J70").Value = 8265
$ws.Range("K70").Value = 5624.353
$ws.Range("L70").Value = 8265
$ws.Range("M70").Value = -5354.353
$ws.Range("N70").Value = -8805

$ws.Range("H73").Value = 6224.5
$ws.Range("I73").Value = 5624.353
$ws.Range("J73").Value = 8265
$ws.Range("K73").Value = 5624.353
$ws.Range("L73").Value = 8265
$ws.Range("M73").Value = -4688.353
$ws.Range("N73").Value = -10137

$ws.Range("H80").Value = 22729628
$ws.Range("I80").Value = 35716500
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 35716500
$ws.Range("L80").Value = 2600
$ws.Range("M80").Value = -35715502
$ws.Range("N80").Value = -4596

$ws.Range("H83").Value = 22729628
$ws.Range("I83").Value = 35716500
$ws.Range("J83").Value = 2600
$ws.Range("K83").Value = 178582500
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = -178577508
$ws.Range("N83").Value = -22984

$ws.Range("H133").Value = 48322.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 48322.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 48322.25
$ws.Range("N133").Value = -58442.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6054.4546
$ws.Range("I7").Value = 2926.25
$ws.Range("J7").Value = 7842
$ws.Range("K7").Value = 2926.25
$ws.Range("L7").Value = 7842
$ws.Range("M7").Value = -2814.25
$ws.Range("N7").Value = -8066

$ws.Range("H61").Value = 2743.5715
$ws.Range("I61").Value = 2800
$ws.Range("J61").Value = 2668.3333
$ws.Range("K61").Value = 2800
$ws.Range("L61").Value = 2668.3333
$ws.Range("M61").Value = -2598
$ws.Range("N61").Value = -3072.3333

$ws.Range("H82").Value = 1527.2941
$ws.Range("I82").Value = 797.4286
$ws.Range("J82").Value = 4933.3335
$ws.Range("K82").Value = 797.4286
$ws.Range("L82").Value = 4933.3335
$ws.Range("M82").Value = -436.4286
$ws.Range("N82").Value = -5655.3335

$ws.Range("H85").Value = 1527.2941
$ws.Range("I85").Value = 797.4286
$ws.Range("J85").Value = 4933.3335
$ws.Range("K85").Value = 797.4286
$ws.Range("L85").Value = 4933.3335
$ws.Range("M85").Value = 450.5714
$ws.Range("N85").Value = -7429.3335

$ws.Range("H113").Value = 2743.5715
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 2668.3333
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 2668.3333
$ws.Range("M113").Value = -630
$ws.Range("N113").Value = -7008.3333

$ws.Range("H122").Value = 3492.6428
$ws.Range("I122").Value = 2554.2727
$ws.Range("J122").Value = 6933.3335
$ws.Range("K122").Value = 7662.8181
$ws.Range("L122").Value = 20800.0005
$ws.Range("M122").Value = -5212.8181
$ws.Range("N122").Value = -25700.0005

$ws.Range("H126").Value = 6054.4546
$ws.Range("I126").Value = 2926.25
$ws.Range("J126").Value = 7842
$ws.Range("K126").Value = 8778.75
$ws.Range("L126").Value = 23526
$ws.Range("M126").Value = -6308.75
$ws.Range("N126").Value = -28466

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null

$ws.Range("H125").Value = 40346.875
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 40346.875
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 40346.875
$ws.Range("N125").Value = -50186.875
